$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-15 07:54:39"
$wsZhCn.Range("G3").Value = "2016-01-15 07:55:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-15 07:54:50"
$wsDeDe.Range("G3").Value = "2016-01-15 07:55:42"
